# Trade #22 closed at 2026-02-17 08:22:10 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: refresh aggregate stats after the new closed trade
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.16   # Current Capital
$summary.Range("B4").Value = -0.84    # Total P&L $
$summary.Range("B5").Value = -0.76    # Total P&L %
$summary.Range("B6").Value = 22       # Total Trades
$summary.Range("B8").Value = 13       # Losing Trades
$summary.Range("B9").Value = 22.73    # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) reflects the same trade
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.16     # Capital
$status.Range("D4").Value = 22        # Trades
$status.Range("E4").Value = -0.84     # P&L $
$status.Range("F4").Value = -0.84     # P&L %
$status.Range("G4").Value = 22.73     # Win Rate %

# ---------------------------------------------------------------------------
# All Trades + MarketMaking sheets: append the new trade row (#22 -> row 23)
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 23

    $ws.Cells.Item($row, 1).Value = 22
    # A leading "'" stops this engine's COM layer from auto-parsing the
    # text as a date (same as real Excel); resetting the style afterwards
    # drops the quote-prefix formatting flag so the cell is plain text,
    # matching the original inline-string cells written by the exporter.
    $ws.Cells.Item($row, 2).Value = "'2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"
    $ws.Cells.Item($row, 3).Value = "08:22:04"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.6899999999999999
    $ws.Cells.Item($row, 7).Value = 0.59
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -14.4928
    $ws.Cells.Item($row, 10).Value = -0.1
    $ws.Cells.Item($row, 11).Value = 99.16
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}
